$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at row 58. We need to add 6 new rows (59-64)
# describing the new "mobility" data fields (Retail/Recreation, Grocery/
# Pharmacy, Parks, Transit, Workplace, Residential), each sharing the same
# look & feel (styles) as the last existing data row (row 58: Category =
# Socioeconomic Data / Indonesia). Copying that row down preserves the
# per-cell formatting (column C & E use cell style index 4, column D uses
# the wrap-text style) exactly like the source workbook does, then we
# overwrite the copied values with the real content.

$ws.Rows.Item(58).Copy()
$ws.Rows.Item(59).Insert(-4121)
$ws.Rows.Item(58).Copy()
$ws.Rows.Item(60).Insert(-4121)
$ws.Rows.Item(58).Copy()
$ws.Rows.Item(61).Insert(-4121)
$ws.Rows.Item(58).Copy()
$ws.Rows.Item(62).Insert(-4121)
$ws.Rows.Item(58).Copy()
$ws.Rows.Item(63).Insert(-4121)
$ws.Rows.Item(58).Copy()
$ws.Rows.Item(64).Insert(-4121)

$newRows = @(
    @("Retail and Recreation Mobility", "csv", "Environmental Data", "Percent change from baseline for mobility in retail and recreation", "Indonesia"),
    @("Grocery and Pharmacy Mobility", "csv", "Environmental Data", "Percent change from baseline for mobility in Grocery and Pharmacy", "Indonesia"),
    @("Parks Mobility", "csv", "Environmental Data", "Percent change from baseline for mobility in parks", "Indonesia"),
    @("Transit Mobility", "csv", "Environmental Data", "Percent change from baseline for mobility for transit", "Indonesia"),
    @("Workplace Mobility", "csv", "Environmental Data", "Percent change from baseline for mobility in workplaces", "Indonesia"),
    @("Residential Mobility", "csv", "Environmental Data", "Percent change from baseline for mobility in residences", "Indonesia")
)

$rowNum = 59
foreach ($rowData in $newRows) {
    $ws.Cells.Item($rowNum, 1).Value = $rowData[0]
    $ws.Cells.Item($rowNum, 2).Value = $rowData[1]
    $ws.Cells.Item($rowNum, 3).Value = $rowData[2]
    $ws.Cells.Item($rowNum, 4).Value = $rowData[3]
    $ws.Cells.Item($rowNum, 5).Value = $rowData[4]
    $rowNum++
}

# The two longest descriptions wrap onto a second line at the column D
# width in use, so those rows end up taller - match that.
$ws.Rows.Item(59).RowHeight = 25.5
$ws.Rows.Item(60).RowHeight = 25.5

# Column A needs to widen to fit the new (longer) header names.
$ws.Columns.Item(1).ColumnWidth = 32.14

# Leave the view scrolled down to the newly-added data, with D64 selected.
[void]$ws.Activate()
[void]$ws.Range("D64").Select()
$excel.ActiveWindow.ScrollRow = 54
$excel.ActiveWindow.ScrollColumn = 1
